$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $text) {
    # Preserve the cell's existing style, force a Text number
    # format so Excel stores the assigned value verbatim (no numeric
    # auto-conversion / float rounding), then restore the original style.
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = $origStyle
}

Set-TextValue $ws.Range("D2") "34.738.93"
$ws.Range("E2").Value = "  +2.62%  "
Set-TextValue $ws.Range("D3") "1.797.00"
$ws.Range("E3").Value = "  -0.37%  "
Set-TextValue $ws.Range("D4") "0.998"
$ws.Range("E4").Value = "  -0.79%  "
Set-TextValue $ws.Range("D5") "225.89"
$ws.Range("E5").Value = "  -0.85%  "
Set-TextValue $ws.Range("D6") "0.556"
$ws.Range("E6").Value = "  -1.79%  "
Set-TextValue $ws.Range("D7") "0.998"
$ws.Range("E7").Value = "  -0.62%  "
Set-TextValue $ws.Range("D8") "32.89"
$ws.Range("E8").Value = "  +4.12%  "
$ws.Range("E9").Value = "  +0.82%  "
Set-TextValue $ws.Range("D10") "0.0678"
$ws.Range("E10").Value = "  +1.69%  "
Set-TextValue $ws.Range("D11") "0.0938"
$ws.Range("E11").Value = "  +0.60%  "
Set-TextValue $ws.Range("D12") "2.054.41"
$ws.Range("E12").Value = "  -0.90%  "
Set-TextValue $ws.Range("D13") "11.21"
$ws.Range("E13").Value = "  +11.94%  "
Set-TextValue $ws.Range("D14") "1.806.38"
$ws.Range("E14").Value = "  +0.08%  "
Set-TextValue $ws.Range("D15") "0.641"
$ws.Range("E15").Value = "  +0.46%  "
Set-TextValue $ws.Range("D16") "34.642.35"
$ws.Range("E16").Value = "  +2.04%  "
Set-TextValue $ws.Range("D17") "4.31"
$ws.Range("E17").Value = "  +2.08%  "
Set-TextValue $ws.Range("D18") "69.87"
$ws.Range("E18").Value = "  +0.30%  "
Set-TextValue $ws.Range("D19") "257.85"
$ws.Range("E19").Value = "  +0.89%  "
Set-TextValue $ws.Range("D20") "0.0₃0780"
$ws.Range("E20").Value = "  +4.40%  "
$ws.Range("E21").Value = "  -0.44%  "
Set-TextValue $ws.Range("D22") "10.49"
$ws.Range("E22").Value = "  +0.17%  "
Set-TextValue $ws.Range("D23") "4.25"
$ws.Range("E23").Value = "  -0.52%  "
Set-TextValue $ws.Range("D24") "2.14"
$ws.Range("E24").Value = "  -2.48%  "
Set-TextValue $ws.Range("D25") "158.65"
$ws.Range("E25").Value = "  -0.54%  "
Set-TextValue $ws.Range("D26") "16.58"
$ws.Range("E26").Value = "  +0.13%  "
Set-TextValue $ws.Range("D27") "7.16"
$ws.Range("E27").Value = "  +2.65%  "
Set-TextValue $ws.Range("D28") "0.115"
$ws.Range("E28").Value = "  -1.21%  "
$ws.Range("E29").Value = "  -0.66%  "
Set-TextValue $ws.Range("D30") "3.82"
$ws.Range("E30").Value = "  -0.62%  "
Set-TextValue $ws.Range("D31") "0.0522"
$ws.Range("E31").Value = "  +0.58%  "
$ws.Range("E32").Value = "  -0.86%  "
Set-TextValue $ws.Range("D33") "3.62"
$ws.Range("E33").Value = "  +1.74%  "
$ws.Range("E34").Value = "  +7.75%  "
Set-TextValue $ws.Range("D35") "1.459.19"
$ws.Range("E35").Value = "  -2.92%  "
$ws.Range("E36").Value = "  -1.25%  "
Set-TextValue $ws.Range("D37") "0.639"
$ws.Range("E37").Value = "  +2.45%  "
Set-TextValue $ws.Range("D38") "0.0191"
$ws.Range("E38").Value = "  +1.49%  "
Set-TextValue $ws.Range("D39") "83.30"
$ws.Range("E39").Value = "  +0.67%  "
$ws.Range("E40").Value = "  +2.65%  "
Set-TextValue $ws.Range("D41") "2.34"
$ws.Range("E41").Value = "  -1.81%  "
Set-TextValue $ws.Range("D42") "0.907"
$ws.Range("E42").Value = "  +0.99%  "
$ws.Range("E43").Value = "  +0.18%  "
Set-TextValue $ws.Range("D44") "0.0510"
$ws.Range("E44").Value = "  -1.94%  "
Set-TextValue $ws.Range("D45") "6.01"
$ws.Range("E45").Value = "  +3.57%  "
$ws.Range("E46").Value = "  -3.06%  "
Set-TextValue $ws.Range("D47") "1.953.98"
$ws.Range("E47").Value = "  -0.76%  "
Set-TextValue $ws.Range("D48") "12.05"
$ws.Range("E48").Value = "  +0.92%  "
$ws.Range("E49").Value = "  -0.47%  "
Set-TextValue $ws.Range("D50") "101.98"
$ws.Range("E50").Value = "  +2.78%  "
Set-TextValue $ws.Range("D51") "49.67"
$ws.Range("E51").Value = "  -2.84%  "
